$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D4").Value = -7.563000000000001
$ws.Range("D7").Value = -7.749200000000009
$ws.Range("D16").Value = -8.212
$ws.Range("D28").Value = -8.196499999999997
$ws.Range("D29").Value = -7.188800000000003
$ws.Range("D32").Value = -7.093999999999991
$ws.Range("D40").Value = -8.842899999999988
$ws.Range("D52").Value = -7.955800000000003
$ws.Range("D57").Value = -8.142099999999994
$ws.Range("D66").Value = -6.918200000000002
$ws.Range("D100").Value = -8.027200000000001
